$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Schedule (Template View")

# Milestone text: "Streamlit/Plotly Dashboard (skills/region/salary)" -> "Dashboard (skills/region/salary)"
$ws.Range("A15").Value = "Dashboard (skills/region/salary)"

# Dates updated - rows 18, 19, 20 all collapse onto 2025-12-08 (serial 45999)
$ws.Range("D18").Value = 45999
$ws.Range("C19").Value = 45999
$ws.Range("D19").Value = 45999
$ws.Range("C20").Value = 45999
$ws.Range("D20").Value = 45999

# Move the selection in the frozen pane
$ws.Range("A26").Select()
